$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh of the cryptos price/volume snapshot (GitHub Actions update run).
# Values that are plain decimal numbers (e.g. "211.92") are written with a
# leading apostrophe so Excel keeps them as text, matching the original
# inline-string cell type instead of silently converting them to numbers.
$ws.Range('D2').Value = '27.923.74'
$ws.Range('E2').Value = '  -0.11%  '
$ws.Range('D3').Value = '1.632.61'
$ws.Range('E3').Value = '  -0.86%  '
$ws.Range('D5').Value = "'211.92"
$ws.Range('E5').Value = '  -0.72%  '
$ws.Range('E6').Value = '  -0.82%  '
$ws.Range('E7').Value = '  -0.11%  '
$ws.Range('D8').Value = "'23.30"
$ws.Range('E8').Value = '  -0.47%  '
$ws.Range('E9').Value = '  -2.78%  '
$ws.Range('E10').Value = '  -0.14%  '
$ws.Range('E11').Value = '  +1.16%  '
$ws.Range('D12').Value = '1.864.67'
$ws.Range('E12').Value = '  -0.84%  '
$ws.Range('D13').Value = '1.631.92'
$ws.Range('E13').Value = '  -0.86%  '
$ws.Range('E14').Value = '  -0.37%  '
$ws.Range('D15').Value = "'0.566"
$ws.Range('E15').Value = '  +0.28%  '
$ws.Range('D16').Value = "'65.28"
$ws.Range('E16').Value = '  -0.44%  '
$ws.Range('D17').Value = '27.924.25'
$ws.Range('E17').Value = '  -0.19%  '
$ws.Range('D18').Value = "'230.58"
$ws.Range('E18').Value = '  -0.88%  '
$ws.Range('E19').Value = '  +0.03%  '
$ws.Range('D20').Value = "'7.52"
$ws.Range('E20').Value = '  -2.20%  '
$ws.Range('E21').Value = '  -0.02%  '
$ws.Range('B22').Value = 'Avalanche'
$ws.Range('C22').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D22').Value = "'10.39"
$ws.Range('E22').Value = '  -2.53%  '
$ws.Range('B23').Value = 'Uniswap'
$ws.Range('C23').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D23').Value = "'4.37"
$ws.Range('E23').Value = '  -0.63%  '
$ws.Range('E24').Value = '  -3.78%  '
$ws.Range('D25').Value = "'154.19"
$ws.Range('E25').Value = '  +1.24%  '
$ws.Range('E26').Value = '  +0.83%  '
$ws.Range('D27').Value = "'15.66"
$ws.Range('E27').Value = '  -0.61%  '
$ws.Range('E28').Value = '  -0.59%  '
$ws.Range('E29').Value = '  -0.10%  '
$ws.Range('E30').Value = '  -0.97%  '
$ws.Range('D31').Value = "'0.0482"
$ws.Range('E31').Value = '  -0.07%  '
$ws.Range('E32').Value = '  +0.96%  '
$ws.Range('D33').Value = "'3.08"
$ws.Range('E33').Value = '  -0.20%  '
$ws.Range('D34').Value = '1.401.90'
$ws.Range('E34').Value = '  -2.88%  '
$ws.Range('E35').Value = '  +0.06%  '
$ws.Range('D36').Value = "'1.00"
$ws.Range('E36').Value = '  +9.10%  '
$ws.Range('E37').Value = '  +1.46%  '
$ws.Range('E38').Value = '  +0.48%  '
$ws.Range('E39').Value = '  +0.49%  '
$ws.Range('D40').Value = "'0.872"
$ws.Range('E40').Value = '  -1.89%  '
$ws.Range('E41').Value = '  -0.20%  '
$ws.Range('D43').Value = "'66.86"
$ws.Range('E43').Value = '  -3.64%  '
$ws.Range('D44').Value = "'5.55"
$ws.Range('E44').Value = '  +2.76%  '
$ws.Range('E45').Value = '  +1.12%  '
$ws.Range('E46').Value = '  -1.26%  '
$ws.Range('D47').Value = '1.774.13'
$ws.Range('E47').Value = '  -0.91%  '
$ws.Range('D48').Value = "'87.73"
$ws.Range('E48').Value = '  -1.39%  '
$ws.Range('E49').Value = '  +0.56%  '
$ws.Range('E50').Value = '  -0.83%  '
$ws.Range('E51').Value = '  -0.27%  '
